$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("REMISIONES   MARZO   2022  ")

# Row 14: fill in payment date (F) and amount collected (G)
$ws.Range("F14").Value = 44638
$ws.Range("G14").Value = 3333

# Row 15: fill in payment date (F) and amount collected (G)
$ws.Range("F15").Value = 44639
$ws.Range("G15").Value = 337

# Row 16: new credit entry
$ws.Range("A16").Value = 44637
$ws.Range("D16").Value = "GUSTAVO"
$ws.Range("E16").Value = 2600
$ws.Range("F16").Value = 44639
$ws.Range("G16").Value = 2600

# Row 17: new credit entry
$ws.Range("A17").Value = 44638
$ws.Range("D17").Value = "MAURO"
$ws.Range("E17").Value = 1080
$ws.Range("F17").Value = 44640
$ws.Range("G17").Value = 1080

# Row 18: new credit entry
$ws.Range("A18").Value = 44639
$ws.Range("D18").Value = "OBRADOR"
$ws.Range("E18").Value = 316
$ws.Range("F18").Value = 44639
$ws.Range("G18").Value = 316

# Row 19: new credit entry
$ws.Range("A19").Value = 44640
$ws.Range("D19").Value = "HERRADURA DAVID"
$ws.Range("E19").Value = 7994
$ws.Range("F19").Value = 44641
$ws.Range("G19").Value = 7994

# Row 20: new credit entry
$ws.Range("A20").Value = 44640
$ws.Range("D20").Value = "HERRADURA DAVID"
$ws.Range("E20").Value = 5256
$ws.Range("F20").Value = 44641
$ws.Range("G20").Value = 5256

# Move the active selection, matching the closing snapshot
$ws.Range("G21").Select()
